$d = $word.ActiveDocument

# Locate the paragraph that discusses primary_contact handling (it is the
# only paragraph containing this literal, ASCII, token).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*primary_contact*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not locate the primary_contact paragraph"
}

$r = $targetPara.Range
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Design choice for handling the primary contact</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The design adopts a dedicated </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>primary_contact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> table to represent the primary contact of a partner. This approach provides a clear and normalized representation of the relationship while enforcing all business rules directly at the database level. It ensures that each partner has at most one primary contact and that the primary contact must belong to the same partner, without introducing circular foreign key dependencies or relying on procedural logic such as triggers. Furthermore, this design preserves Third Normal Form (3NF) and Boyce–Codd Normal Form (BCNF) by separating the primary-contact relationship from the partner and contact entities, thereby avoiding redundancy and update anomalies.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Other design alternatives were considered but not selected. Storing the primary contact directly in the </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>partner</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> table was rejected due to circular foreign key dependencies, which complicate data insertion, deletion, and maintenance. Modeling the primary contact as a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> attribute of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>contact_point</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> entity was not chosen because it embeds relationship semantics into an entity attribute and requires triggers, leading to hidden constraints and reduced maintainability. A more generalized design using </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2 </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>junction tables</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>partner_primary_contact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>and</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t>partner_contact</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>was also considered but deemed unnecessarily complex for the current requirements, as contacts are not expected to be shared across multiple partners.</w:t></w:r></w:p>'
$r.InsertXML($xmlFrag)
